# NATMI LR-pair output (Cthrc1 -> Fzd6) recomputed after the upstream TPM matrix was updated.
# Only the 'ECs' cluster's Cthrc1 (ligand) / Fzd6 (receptor) expression changed; every derived
# specificity/weight column that is normalised across sending or target clusters is re-stamped
# for all 9 Sending x Target rows as a consequence.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026182
$ws.Range("H2").Value = 0.078546
$ws.Range("I2").Value = 0.02060098669457318
$ws.Range("J2").Value = 0.02060098669457318
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.28101533333333
$ws.Range("N2").Value = 36.843046
$ws.Range("O2").Value = 0.959552102275422
$ws.Range("P2").Value = 0.959552102275422
$ws.Range("Q2").Value = 0.3215415434573334
$ws.Range("R2").Value = 2.893873891116
$ws.Range("S2").Value = 0.01976772009172569
$ws.Range("T2").Value = 0.01976772009172569

# Row 3: Sending=ECs, Target=FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.026182
$ws.Range("H3").Value = 0.078546
$ws.Range("I3").Value = 0.02060098669457318
$ws.Range("J3").Value = 0.02060098669457318
$ws.Range("M3").Value = 0.09168666666666665
$ws.Range("O3").Value = 0.007163750827004844
$ws.Range("P3").Value = 0.007163750827004845
$ws.Range("Q3").Value = 0.002400540306666667
$ws.Range("R3").Value = 0.02160486276
$ws.Range("S3").Value = 0.0001475803354703644
$ws.Range("T3").Value = 0.0001475803354703644

# Row 4: Sending=ECs, Target=MuSCs
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.026182
$ws.Range("H4").Value = 0.078546
$ws.Range("I4").Value = 0.02060098669457318
$ws.Range("J4").Value = 0.02060098669457318
$ws.Range("M4").Value = 0.4259936666666666
$ws.Range("O4").Value = 0.03328414689757318
$ws.Range("P4").Value = 0.03328414689757318
$ws.Range("Q4").Value = 0.01115336618066667
$ws.Range("R4").Value = 0.100380295626
$ws.Range("S4").Value = 0.0006856862673771243
$ws.Range("T4").Value = 0.0006856862673771242

# Row 5: Sending=FAPs, Target=ECs
$ws.Range("I5").Value = 0.9231010325934437
$ws.Range("J5").Value = 0.9231010325934434
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.28101533333333
$ws.Range("N5").Value = 36.843046
$ws.Range("O5").Value = 0.959552102275422
$ws.Range("P5").Value = 0.959552102275422
$ws.Range("Q5").Value = 14.40782110040111
$ws.Range("R5").Value = 129.67038990361
$ws.Range("S5").Value = 0.8857635364376517
$ws.Range("T5").Value = 0.8857635364376515

# Row 6: Sending=FAPs, Target=FAPs
$ws.Range("I6").Value = 0.9231010325934437
$ws.Range("J6").Value = 0.9231010325934434
$ws.Range("M6").Value = 0.09168666666666665
$ws.Range("O6").Value = 0.007163750827004844
$ws.Range("P6").Value = 0.007163750827004845
$ws.Range("R6").Value = 0.9680832971
$ws.Range("S6").Value = 0.006612865785650308
$ws.Range("T6").Value = 0.006612865785650307

# Row 7: Sending=FAPs, Target=MuSCs
$ws.Range("I7").Value = 0.9231010325934437
$ws.Range("J7").Value = 0.9231010325934434
$ws.Range("M7").Value = 0.4259936666666666
$ws.Range("O7").Value = 0.03328414689757318
$ws.Range("P7").Value = 0.03328414689757318
$ws.Range("Q7").Value = 0.4997665398705555
$ws.Range("R7").Value = 4.497898858835
$ws.Range("S7").Value = 0.03072463037014166
$ws.Range("T7").Value = 0.03072463037014166

# Row 8: Sending=MuSCs, Target=ECs
$ws.Range("I8").Value = 0.05629798071198328
$ws.Range("J8").Value = 0.05629798071198327
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.28101533333333
$ws.Range("N8").Value = 36.843046
$ws.Range("O8").Value = 0.959552102275422
$ws.Range("P8").Value = 0.959552102275422
$ws.Range("Q8").Value = 0.8787025534282222
$ws.Range("R8").Value = 7.908322980854001
$ws.Range("S8").Value = 0.05402084574604471
$ws.Range("T8").Value = 0.0540208457460447

# Row 9: Sending=MuSCs, Target=FAPs
$ws.Range("I9").Value = 0.05629798071198328
$ws.Range("J9").Value = 0.05629798071198327
$ws.Range("M9").Value = 0.09168666666666665
$ws.Range("O9").Value = 0.007163750827004844
$ws.Range("P9").Value = 0.007163750827004845
$ws.Range("Q9").Value = 0.006560150437777777
$ws.Range("R9").Value = 0.05904135393999999
$ws.Range("S9").Value = 0.0004033047058841729
$ws.Range("T9").Value = 0.0004033047058841729

# Row 10: Sending=MuSCs, Target=MuSCs
$ws.Range("I10").Value = 0.05629798071198328
$ws.Range("J10").Value = 0.05629798071198327
$ws.Range("M10").Value = 0.4259936666666666
$ws.Range("O10").Value = 0.03328414689757318
$ws.Range("P10").Value = 0.03328414689757318
$ws.Range("S10").Value = 0.001873830260054393
$ws.Range("T10").Value = 0.001873830260054392
